$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 54262.6
$ws.Range("J3").Value = 54262.6
$ws.Range("L3").Value = 54262.6
$ws.Range("N3").Value = -54490.6
$ws.Range("H70").Value = 3812.3333
$ws.Range("I70").Value = 3949
$ws.Range("K70").Value = 11847
$ws.Range("M70").Value = -11577
$ws.Range("H73").Value = 3812.3333
$ws.Range("I73").Value = 3949
$ws.Range("K73").Value = 11847
$ws.Range("M73").Value = -10911
$ws.Range("H102").Value = 54262.6
$ws.Range("J102").Value = 54262.6
$ws.Range("L102").Value = 54262.6
$ws.Range("N102").Value = -60752.6
$ws.Range("H138").Value = 2928
$ws.Range("I138").Value = 2832
$ws.Range("K138").Value = 8496
$ws.Range("M138").Value = -3356

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2665
$ws.Range("I2").Value = 2665
$ws.Range("K2").Value = 2665
$ws.Range("M2").Value = -2552
$ws.Range("H31").Value = 23999.8
$ws.Range("I31").Value = 23999.8
$ws.Range("K31").Value = 23999.8
$ws.Range("M31").Value = -23705.8
$ws.Range("H32").Value = 1910.8572
$ws.Range("I32").Value = 1672.9412
$ws.Range("K32").Value = 1672.9412
$ws.Range("M32").Value = -1385.9412
$ws.Range("H45").Value = 751.75
$ws.Range("I45").Value = 751.75
$ws.Range("K45").Value = 751.75
$ws.Range("M45").Value = -374.75
$ws.Range("H61").Value = 6740.25
$ws.Range("I61").Value = 3480.75
$ws.Range("J61").Value = 9999.75
$ws.Range("K61").Value = 3480.75
$ws.Range("L61").Value = 9999.75
$ws.Range("M61").Value = -3268.75
$ws.Range("N61").Value = -10423.75
$ws.Range("H74").Value = 3845
$ws.Range("I74").Value = 3727.75
$ws.Range("J74").Value = 4314
$ws.Range("K74").Value = 3727.75
$ws.Range("L74").Value = 4314
$ws.Range("M74").Value = -2853.75
$ws.Range("N74").Value = -6062
$ws.Range("H77").Value = 3845
$ws.Range("I77").Value = 3727.75
$ws.Range("J77").Value = 4314
$ws.Range("K77").Value = 18638.75
$ws.Range("L77").Value = 21570
$ws.Range("M77").Value = -14270.75
$ws.Range("N77").Value = -30306
$ws.Range("H102").Value = 2100
$ws.Range("I102").Value = 2100
$ws.Range("K102").Value = 2100
$ws.Range("M102").Value = -478
$ws.Range("H110").Value = 3480.5454
$ws.Range("I110").Value = 3373.75
$ws.Range("J110").Value = 3541.5715
$ws.Range("K110").Value = 3373.75
$ws.Range("L110").Value = 3541.5715
$ws.Range("M110").Value = -1328.75
$ws.Range("N110").Value = -7631.5715
$ws.Range("H116").Value = 2665
$ws.Range("I116").Value = 2665
$ws.Range("K116").Value = 2665
$ws.Range("M116").Value = -371
$ws.Range("H136").Value = 6740.25
$ws.Range("I136").Value = 3480.75
$ws.Range("J136").Value = 9999.75
$ws.Range("K136").Value = 10442.25
$ws.Range("L136").Value = 29999.25
$ws.Range("M136").Value = -7892.25
$ws.Range("N136").Value = -35099.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2665
$ws.Range("I3").Value = 2665
$ws.Range("K3").Value = 2665
$ws.Range("M3").Value = -2551
$ws.Range("H64").Value = 418.6
$ws.Range("I64").Value = 248.33333
$ws.Range("K64").Value = 248.33333
$ws.Range("M64").Value = -23.33332999999999
$ws.Range("H67").Value = 418.6
$ws.Range("I67").Value = 248.33333
$ws.Range("K67").Value = 248.33333
$ws.Range("M67").Value = 531.6666700000001
$ws.Range("H86").Value = 3298
$ws.Range("I86").Value = 3916.6667
$ws.Range("K86").Value = 3916.6667
$ws.Range("M86").Value = -2793.6667
$ws.Range("H89").Value = 3298
$ws.Range("I89").Value = 3916.6667
$ws.Range("K89").Value = 19583.3335
$ws.Range("M89").Value = -13967.3335
$ws.Range("H94").Value = 598.4286
$ws.Range("I94").Value = 598.4286
$ws.Range("K94").Value = 598.4286
$ws.Range("M94").Value = -147.4286
$ws.Range("H95").Value = 28556.334
$ws.Range("J95").Value = 28556.334
$ws.Range("L95").Value = 28556.334
$ws.Range("N95").Value = -34048.334
$ws.Range("H105").Value = 1856
$ws.Range("I105").Value = 1831.3334
$ws.Range("K105").Value = 1831.3334
$ws.Range("M105").Value = -84.33339999999998
$ws.Range("H106").Value = 7035.5
$ws.Range("J106").Value = 7035.5
$ws.Range("L106").Value = 7035.5
$ws.Range("N106").Value = -9559.5
$ws.Range("H134").Value = 4002.5186
$ws.Range("I134").Value = 4110.3076
$ws.Range("J134").Value = 1200
$ws.Range("K134").Value = 12330.9228
$ws.Range("L134").Value = 3600
$ws.Range("M134").Value = -9795.9228
$ws.Range("N134").Value = -8670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2636
$ws.Range("I31").Value = 2636
$ws.Range("K31").Value = 2636
$ws.Range("M31").Value = -2341
$ws.Range("H34").Value = 2636
$ws.Range("I34").Value = 2636
$ws.Range("K34").Value = 2636
$ws.Range("M34").Value = -2434
$ws.Range("H43").Value = 8997.6
$ws.Range("J43").Value = 8997.6
$ws.Range("L43").Value = 8997.6
$ws.Range("N43").Value = -9365.6
$ws.Range("H101").Value = 8997.6
$ws.Range("J101").Value = 8997.6
$ws.Range("L101").Value = 8997.6
$ws.Range("N101").Value = -15487.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 51.615383
$ws.Range("I2").Value = 64.44444
$ws.Range("J2").Value = 22.75
$ws.Range("K2").Value = 386.66664
$ws.Range("L2").Value = 136.5
$ws.Range("M2").Value = -273.66664
$ws.Range("N2").Value = -362.5
$ws.Range("H4").Value = 2106.6365
$ws.Range("I4").Value = 2204
$ws.Range("K4").Value = 6612
$ws.Range("M4").Value = -6500
$ws.Range("H10").Value = 35.625
$ws.Range("I10").Value = 26.428572
$ws.Range("K10").Value = 79.28571599999999
$ws.Range("M10").Value = 59.71428400000001
$ws.Range("H14").Value = 25292.75
$ws.Range("I14").Value = 25292.75
$ws.Range("K14").Value = 75878.25
$ws.Range("M14").Value = -75705.25
$ws.Range("H25").Value = 1385.8334
$ws.Range("I25").Value = 945
$ws.Range("K25").Value = 2835
$ws.Range("M25").Value = -2666
$ws.Range("H30").Value = 1385.8334
$ws.Range("I30").Value = 945
$ws.Range("K30").Value = 2835
$ws.Range("M30").Value = -2733
$ws.Range("H36").Value = 175
$ws.Range("I36").Value = 175
$ws.Range("K36").Value = 525
$ws.Range("M36").Value = -356
$ws.Range("H121").Value = 831.6
$ws.Range("I121").Value = 789.5
$ws.Range("K121").Value = 2368.5
$ws.Range("M121").Value = -1058.5
$ws.Range("H139").Value = 1392.5
$ws.Range("I139").Value = 1392.5
$ws.Range("K139").Value = 4177.5
$ws.Range("M139").Value = 962.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 2680.111
$ws.Range("I9").Value = 531.5714
$ws.Range("K9").Value = 531.5714
$ws.Range("M9").Value = -361.5714
$ws.Range("H105").Value = 217500
$ws.Range("J105").Value = 217500
$ws.Range("L105").Value = 217500
$ws.Range("N105").Value = -224488
$ws.Range("H123").Value = 26666.666
$ws.Range("J123").Value = 26666.666
$ws.Range("L123").Value = 26666.666
$ws.Range("N123").Value = -31566.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1233.4
$ws.Range("I9").Value = 1325
$ws.Range("K9").Value = 1325
$ws.Range("M9").Value = -1101
$ws.Range("H97").Value = 200000
$ws.Range("J97").Value = 200000
$ws.Range("L97").Value = 200000
$ws.Range("N97").Value = -201982
$ws.Range("H122").Value = 3878
$ws.Range("J122").Value = 5000
$ws.Range("L122").Value = 15000
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 54635
$ws.Range("J68").Value = 54635
$ws.Range("L68").Value = 54635
$ws.Range("N68").Value = -56257
$ws.Range("H71").Value = 54635
$ws.Range("J71").Value = 54635
$ws.Range("L71").Value = 163905
$ws.Range("N71").Value = -172017
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H105").Value = 115000
$ws.Range("J105").Value = 115000
$ws.Range("L105").Value = 115000
$ws.Range("N105").Value = -121988
$ws.Range("H132").Value = 1923.25
$ws.Range("I132").Value = 1651.4667
$ws.Range("K132").Value = 4954.4001
$ws.Range("M132").Value = -2424.4001
